# Weekly price update for "Apio" (celery) sheet.
# A new week's data (2 rows: Primera + Segunda quality) is inserted right
# before the existing row 1037, pushing all subsequent rows down by two
# positions (old 1037 -> new 1039, ... old 1068 -> new 1070).
# The dimension grows from A1:R1068 to A1:R1070 automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 1037, shifting everything below (old rows
# 1037..1068) down to 1039..1070.
$ws.Rows("1037:1038").Insert()

# New row 1037: Primera quality, week of 2023-05-29 (serial 45075).
$ws.Range("A1037").Value = 6
$ws.Range("B1037").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C1037").Value = 'Metropolitana'
$ws.Range("D1037").Value = 45075
$ws.Range("E1037").Value = 13
$ws.Range("F1037").Value = 100112017
$ws.Range("G1037").Value = 'Apio'
$ws.Range("H1037").Value = 'Americana (o)'
$ws.Range("I1037").Value = 'Primera'
$ws.Range("J1037").Value = 2000
$ws.Range("K1037").Value = 6000
$ws.Range("L1037").Value = 7000
$ws.Range("M1037").Value = 6600
$ws.Range("N1037").Value = '$/docena de matas'
$ws.Range("O1037").Value = 'Región de Coquimbo'
$ws.Range("P1037").Value = 1100
$ws.Range("Q1037").Value = 6
$ws.Range("R1037").Value = 'Hortaliza'

# New row 1038: Segunda quality, same week.
$ws.Range("A1038").Value = 6
$ws.Range("B1038").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C1038").Value = 'Metropolitana'
$ws.Range("D1038").Value = 45075
$ws.Range("E1038").Value = 13
$ws.Range("F1038").Value = 100112017
$ws.Range("G1038").Value = 'Apio'
$ws.Range("H1038").Value = 'Americana (o)'
$ws.Range("I1038").Value = 'Segunda'
$ws.Range("J1038").Value = 900
$ws.Range("K1038").Value = 4000
$ws.Range("L1038").Value = 4000
$ws.Range("M1038").Value = 4000
$ws.Range("N1038").Value = '$/docena de matas'
$ws.Range("O1038").Value = 'Región de Coquimbo'
$ws.Range("P1038").Value = 667
$ws.Range("Q1038").Value = 6
$ws.Range("R1038").Value = 'Hortaliza'
